# Initial Averaging for Clinical Epitope Ranking done
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: F6:P6 -> 7 (fills in previously blank L6/M6 as well)
$ws.Range("F6:P6").Value = 7

# Row 7: F7:P7 -> 20 (replaces previous "Error" text in L7 and blank M7)
$ws.Range("F7:P7").Value = 20

# Row 8: F8:P8 -> 30 (replaces previous "Error" text in L8 and blank M8)
$ws.Range("F8:P8").Value = 30

# Update sheet view: scroll so column D is the leftmost visible column,
# and move the active selection to V10
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("V10").Select()
